# [SPARK-2419][Streaming][Docs] More updates to the streaming programming guide
# - Refresh the cached "datetimeFigureOut" date placeholder text (slide master + all
#   slide layouts) from 1/23/14 to 9/6/14.
# - Update the streaming architecture diagram on slide 1: "HDFS" -> "HDFS/S3" and
#   "ZeroMQ" -> "Kinesis".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached date field text ("Date Placeholder ...") everywhere it
#    appears: once on the slide master, and once on every custom (slide)
#    layout attached to that master.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapeRange) {
    for ($i = 1; $i -le $shapeRange.Count; $i++) {
        $shp = $shapeRange.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
                $shp.TextFrame.TextRange.Text = "9/6/14"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Update the streaming sources diagram text on slide 1.
#    The diagram is a deeply nested group; GroupItems() on the outermost
#    group flattens to every leaf shape, so we can find our targets by name.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$outerGroup = $slide1.Shapes.Item(1)
$items = $outerGroup.GroupItems

for ($i = 1; $i -le $items.Count; $i++) {
    $shp = $items.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }

    $text = $shp.TextFrame.TextRange.Text
    if ($shp.Name -eq "Rounded Rectangle 47" -and $text -eq "HDFS") {
        $shp.TextFrame.TextRange.Text = "HDFS/S3"
    }
    elseif ($shp.Name -eq "Rounded Rectangle 48" -and $text -eq "ZeroMQ") {
        $shp.TextFrame.TextRange.Text = "Kinesis"
    }
}
